$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new reporting quarter (period ending 2018-09-30) was added to the financials.
# Insert a fresh column at D; this pushes the existing D:K data one column right (to E:L).
$ws.Columns("D:D").Insert()

# Clone number formatting/styles for the new column from its neighbour (old column D, now E)
# so the new quarter column renders with the same date/number formats as the rest of the table.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Columns("D:D").ColumnWidth = $ws.Columns("E:E").ColumnWidth

# Fill in the new quarter's figures in column D
$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 46800
$ws.Range("D9").Value = 34300
$ws.Range("D10").Value = 12500
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 45300
$ws.Range("D18").Value = 1500
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = 3100
$ws.Range("D22").Value = 200
$ws.Range("D23").Value = 1300
$ws.Range("D24").Value = 500
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 800
$ws.Range("D27").Value = 800
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = 800
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 800
$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 700
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 33900
$ws.Range("D44").Value = 29600
$ws.Range("D45").Value = 2000
$ws.Range("D46").Value = 66200
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 51900
$ws.Range("D49").Value = "NA"
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 24900
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 143000
$ws.Range("D57").Value = 15400
$ws.Range("D58").Value = 13400
$ws.Range("D59").Value = 28000
$ws.Range("D60").Value = 56700
$ws.Range("D61").Value = 22200
$ws.Range("D62").Value = 34500
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 113400
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 62000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 29600
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = 800
$ws.Range("D83").Value = 1600
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 4100
$ws.Range("D91").Value = -1000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -1000
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -2900
$ws.Range("D101").Value = -400
$ws.Range("D102").Value = -100
